$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ------------------------------------------------------------------
# 1) Rename the header row: "<Name>_old" -> "<Name>_FV2310"
#                            "<Name>_new" -> "<Name>_FV2404"
# ------------------------------------------------------------------
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -like "*_old") {
        $cell.Value2 = ($val -replace "_old$", "_FV2310")
    } elseif ($val -like "*_new") {
        $cell.Value2 = ($val -replace "_new$", "_FV2404")
    }
}

# ------------------------------------------------------------------
# 2) Turn the header/data range into a real Excel Table (ListObject)
#    so the headers also show up as table column names.
# ------------------------------------------------------------------
$range = $ws.Range("A1:U85")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"

# ------------------------------------------------------------------
# 3) Freeze the header row (row 1) so it stays visible when scrolling.
# ------------------------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
